$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.136.03'
$ws.Range('E2').Value = '  +0.37%  '
$ws.Range('D3').Value = '2.931.48'
$ws.Range('E3').Value = '  +0.72%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''593.34'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('D6').Value = '''145.11'
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '''0.504'
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('D9').Value = '''6.98'
$ws.Range('E9').Value = '  +4.11%  '
$ws.Range('E10').Value = '  -0.81%  '
$ws.Range('E11').Value = '  -0.55%  '
$ws.Range('E12').Value = '  -0.41%  '
$ws.Range('E13').Value = '  +0.65%  '
$ws.Range('D15').Value = '3.417.71'
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').Value = '61.064.95'
$ws.Range('E16').Value = '  +0.39%  '
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('D18').Value = '2.933.37'
$ws.Range('E18').Value = '  +1.13%  '
$ws.Range('D19').Value = '''433.73'
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('D20').Value = '''13.49'
$ws.Range('E20').Value = '  -0.31%  '
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('D22').Value = '''7.11'
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('D23').Value = '''81.76'
$ws.Range('E23').Value = '  +0.36%  '
$ws.Range('D24').Value = '''11.06'
$ws.Range('E24').Value = '  +2.14%  '
$ws.Range('E25').Value = '  -0.61%  '
$ws.Range('D26').Value = '''11.87'
$ws.Range('E26').Value = '  -1.20%  '
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  -2.07%  '
$ws.Range('D29').Value = '''2.62'
$ws.Range('E29').Value = '  -0.31%  '
$ws.Range('D30').Value = '''6.98'
$ws.Range('E30').Value = '  -1.23%  '
$ws.Range('E31').Value = '  +3.12%  '
$ws.Range('D32').Value = '''26.72'
$ws.Range('E32').Value = '  +0.38%  '
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('D34').Value = '0.0₃0871'
$ws.Range('E34').Value = '  +1.60%  '
$ws.Range('D35').Value = '''1.01'
$ws.Range('E35').Value = '  +0.59%  '
$ws.Range('E36').Value = '  +0.63%  '
$ws.Range('E37').Value = '  -0.98%  '
$ws.Range('E38').Value = '  -1.06%  '
$ws.Range('D39').Value = '''2.00'
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('E40').Value = '  -0.20%  '
$ws.Range('E41').Value = '  +4.72%  '
$ws.Range('E42').Value = '  -2.52%  '
$ws.Range('D43').Value = '''374.56'
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('D44').Value = '''0.0347'
$ws.Range('E44').Value = '  -0.65%  '
$ws.Range('D45').Value = '2.707.71'
$ws.Range('D46').Value = '''133.48'
$ws.Range('E46').Value = '  +2.67%  '
$ws.Range('D48').Value = '''23.93'
$ws.Range('E48').Value = '  -0.62%  '
$ws.Range('E49').Value = '  -1.02%  '
$ws.Range('E50').Value = '  -1.23%  '
$ws.Range('E51').Value = '  +0.10%  '
